# Auto-generated edit script: updates currentAveragePrice/Leve profit
# figures across multiple job sheets, matching the upstream market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1375.28
$ws.Range("I15").Value = 1375.28
$ws.Range("K15").Value = 4125.84
$ws.Range("M15").Value = -3956.84
$ws.Range("H17").Value = 246.08333
$ws.Range("J17").Value = 246.08333
$ws.Range("L17").Value = 738.24999
$ws.Range("N17").Value = -1074.24999
$ws.Range("H86").Value = 7442
$ws.Range("I86").Value = 1720
$ws.Range("J86").Value = 18886
$ws.Range("K86").Value = 1720
$ws.Range("L86").Value = 18886
$ws.Range("M86").Value = -597
$ws.Range("N86").Value = -21132
$ws.Range("H89").Value = 7442
$ws.Range("I89").Value = 1720
$ws.Range("J89").Value = 18886
$ws.Range("K89").Value = 8600
$ws.Range("L89").Value = 94430
$ws.Range("M89").Value = -2984
$ws.Range("N89").Value = -105662
$ws.Range("H92").Value = 515.44446
$ws.Range("I92").Value = 515.44446
$ws.Range("K92").Value = 515.44446
$ws.Range("M92").Value = 732.55554
$ws.Range("H128").Value = 43360.855
$ws.Range("J128").Value = 43360.855
$ws.Range("L128").Value = 43360.855
$ws.Range("N128").Value = -53320.855
$ws.Range("H135").Value = 1038.7693
$ws.Range("I135").Value = 1045.8182
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 9412.363799999999
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -6877.363799999999
$ws.Range("N135").Value = -14070
$ws.Range("H138").Value = 4319.5186
$ws.Range("J138").Value = 5399.7026
$ws.Range("L138").Value = 16199.1078
$ws.Range("N138").Value = -26479.1078

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 6521.875
$ws.Range("J44").Value = 6521.875
$ws.Range("L44").Value = 6521.875
$ws.Range("N44").Value = -7497.875
$ws.Range("H45").Value = 1530.3077
$ws.Range("I45").Value = 1584.2858
$ws.Range("J45").Value = 1467.3334
$ws.Range("K45").Value = 1584.2858
$ws.Range("L45").Value = 1467.3334
$ws.Range("M45").Value = -1207.2858
$ws.Range("N45").Value = -2221.3334
$ws.Range("H74").Value = 2711.0952
$ws.Range("I74").Value = 1297
$ws.Range("J74").Value = 5539.2856
$ws.Range("K74").Value = 1297
$ws.Range("L74").Value = 5539.2856
$ws.Range("M74").Value = -423
$ws.Range("N74").Value = -7287.2856
$ws.Range("H77").Value = 2711.0952
$ws.Range("I77").Value = 1297
$ws.Range("J77").Value = 5539.2856
$ws.Range("K77").Value = 6485
$ws.Range("L77").Value = 27696.428
$ws.Range("M77").Value = -2117
$ws.Range("N77").Value = -36432.428
$ws.Range("H112").Value = 15236.6
$ws.Range("J112").Value = 15236.6
$ws.Range("L112").Value = 15236.6
$ws.Range("N112").Value = -18190.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43506.582
$ws.Range("I20").Value = 68849.2
$ws.Range("J20").Value = 1268.8889
$ws.Range("K20").Value = 68849.2
$ws.Range("L20").Value = 1268.8889
$ws.Range("M20").Value = -68602.2
$ws.Range("N20").Value = -1762.8889
$ws.Range("H35").Value = 18139
$ws.Range("J35").Value = 19766.8
$ws.Range("L35").Value = 19766.8
$ws.Range("N35").Value = -20386.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9760
$ws.Range("H27").Value = 10000
$ws.Range("I27").Value = 10000
$ws.Range("K27").Value = 10000
$ws.Range("M27").Value = -9808
$ws.Range("H122").Value = 4317.7393
$ws.Range("I122").Value = 3826
$ws.Range("J122").Value = 5239.75
$ws.Range("K122").Value = 11478
$ws.Range("L122").Value = 15719.25
$ws.Range("M122").Value = -9028
$ws.Range("N122").Value = -20619.25
$ws.Range("I132").Value = 45458230
$ws.Range("K132").Value = 136374690
$ws.Range("M132").Value = -136372160
$ws.Range("H134").Value = 1143.8667
$ws.Range("I134").Value = 1118.4286
$ws.Range("K134").Value = 3355.2858
$ws.Range("M134").Value = -820.2857999999997

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100001610
$ws.Range("J4").Value = 111112890
$ws.Range("L4").Value = 333338670
$ws.Range("N4").Value = -333338894
$ws.Range("H33").Value = 3467.2222
$ws.Range("I33").Value = 3385
$ws.Range("J33").Value = 3508.3333
$ws.Range("K33").Value = 20310
$ws.Range("L33").Value = 21049.9998
$ws.Range("M33").Value = -20027
$ws.Range("N33").Value = -21615.9998
$ws.Range("H34").Value = 2249.9
$ws.Range("J34").Value = 2774.875
$ws.Range("L34").Value = 8324.625
$ws.Range("N34").Value = -8492.625
$ws.Range("H39").Value = 4250
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4250
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 12750
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -13338
$ws.Range("H55").Value = 14862.294
$ws.Range("I55").Value = 50300
$ws.Range("J55").Value = 10137.267
$ws.Range("K55").Value = 150900
$ws.Range("L55").Value = 30411.801
$ws.Range("M55").Value = -150723
$ws.Range("N55").Value = -30765.801

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 502426.25
$ws.Range("I102").Value = 2889.375
$ws.Range("K102").Value = 2889.375
$ws.Range("M102").Value = -1267.375
$ws.Range("H135").Value = 27815.264
$ws.Range("J135").Value = 27815.264
$ws.Range("L135").Value = 27815.264
$ws.Range("N135").Value = -37955.264

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1877.5358
$ws.Range("I7").Value = 1532.9048
$ws.Range("J7").Value = 2911.4285
$ws.Range("K7").Value = 1532.9048
$ws.Range("L7").Value = 2911.4285
$ws.Range("M7").Value = -1420.9048
$ws.Range("N7").Value = -3135.4285
$ws.Range("H122").Value = 3768
$ws.Range("I122").Value = 4999.5
$ws.Range("J122").Value = 1305
$ws.Range("K122").Value = 14998.5
$ws.Range("L122").Value = 3915
$ws.Range("M122").Value = -12548.5
$ws.Range("N122").Value = -8815
$ws.Range("H126").Value = 1877.5358
$ws.Range("I126").Value = 1532.9048
$ws.Range("J126").Value = 2911.4285
$ws.Range("K126").Value = 4598.7144
$ws.Range("L126").Value = 8734.2855
$ws.Range("M126").Value = -2128.7144
$ws.Range("N126").Value = -13674.2855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 250858.62
$ws.Range("I81").Value = 200590
$ws.Range("J81").Value = 334639.66
$ws.Range("K81").Value = 401180
$ws.Range("L81").Value = 669279.3199999999
$ws.Range("M81").Value = -400119
$ws.Range("N81").Value = -671401.3199999999
$ws.Range("H84").Value = 250858.62
$ws.Range("I84").Value = 200590
$ws.Range("J84").Value = 334639.66
$ws.Range("K84").Value = 2005900
$ws.Range("L84").Value = 3346396.6
$ws.Range("M84").Value = -2000596
$ws.Range("N84").Value = -3357004.6

